# Update the "dSF" (column F) values for the webb_logan 2022 sheet.
# These reflect repulled data / recalculated mean values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = 4
    4  = -4
    5  = 2
    6  = 3
    8  = -3
    9  = -4
    10 = 2
    11 = 2
    12 = -5
    15 = 8
    16 = 8
    17 = -1
    18 = -4
    19 = 5
    20 = -1
    21 = 3
    22 = -1
    23 = 2
    24 = 1
    25 = 4
    26 = 1
    27 = 2
    28 = 5
    29 = 1
    30 = 3
    31 = -3
    32 = 1
    33 = 3
    34 = -3
    35 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
